# Fix typo in predictor labels: "per capita" -> "per cap." (and fix
# mismatched bracket/parenthesis in the livestock consumption label).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "C2"  = "ln(GDP [dollars per cap.])"
    "C3"  = "ln(ProMed Mentions [per cap.])"
    "C4"  = "ln(Migrant Population [per cap.])"
    "C6"  = "ln(Tourism - Inbound [per cap.])"
    "C8"  = "ln(AB Exports [dollars per cap.])"
    "C9"  = "ln(Publication Bias Index [per cap.])"
    "C12" = "Livestock AB Consumption (kg per cap.)"
    "C13" = "ln(ProMed Mentions [per cap.])"
    "C14" = "ln(Publication Bias Index [per cap.])"
    "C17" = "ln(GDP [dollars per cap.])"
}

foreach ($addr in $replacements.Keys) {
    $ws.Range($addr).Value = $replacements[$addr]
}
